$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (reshuffled weekly data) ---

$ws.Cells.Item(2,4).Value = 44377
$ws.Cells.Item(2,11).Value = 11000
$ws.Cells.Item(2,12).Value = 12000
$ws.Cells.Item(2,13).Value = 11500
$ws.Cells.Item(2,16).Value = 192

$ws.Cells.Item(3,4).Value = 44350
$ws.Cells.Item(3,10).Value = 60

$ws.Cells.Item(4,4).Value = 44327
$ws.Cells.Item(4,11).Value = 9000
$ws.Cells.Item(4,12).Value = 10000
$ws.Cells.Item(4,13).Value = 9500
$ws.Cells.Item(4,16).Value = 158

$ws.Cells.Item(5,4).Value = 44383

$ws.Cells.Item(6,4).Value = 44336
$ws.Cells.Item(6,11).Value = 12000
$ws.Cells.Item(6,12).Value = 13000
$ws.Cells.Item(6,13).Value = 12500
$ws.Cells.Item(6,16).Value = 208

$ws.Cells.Item(7,4).Value = 44166
$ws.Cells.Item(7,11).Value = 7000
$ws.Cells.Item(7,12).Value = 7500
$ws.Cells.Item(7,13).Value = 7250
$ws.Cells.Item(7,16).Value = 121

$ws.Cells.Item(8,4).Value = 44308

$ws.Cells.Item(9,4).Value = 44189
$ws.Cells.Item(9,12).Value = 8000
$ws.Cells.Item(9,13).Value = 7500
$ws.Cells.Item(9,16).Value = 125

$ws.Cells.Item(10,4).Value = 44355
$ws.Cells.Item(10,11).Value = 11000
$ws.Cells.Item(10,12).Value = 12000
$ws.Cells.Item(10,13).Value = 11500
$ws.Cells.Item(10,16).Value = 192

$ws.Cells.Item(11,4).Value = 44246

$ws.Cells.Item(12,4).Value = 44230
$ws.Cells.Item(12,11).Value = 8000
$ws.Cells.Item(12,12).Value = 9000
$ws.Cells.Item(12,13).Value = 8500
$ws.Cells.Item(12,16).Value = 142

$ws.Cells.Item(13,4).Value = 44159
$ws.Cells.Item(13,11).Value = 8000
$ws.Cells.Item(13,12).Value = 9000
$ws.Cells.Item(13,13).Value = 8500
$ws.Cells.Item(13,16).Value = 142

$ws.Cells.Item(14,4).Value = 44328
$ws.Cells.Item(14,11).Value = 9000
$ws.Cells.Item(14,12).Value = 10000
$ws.Cells.Item(14,13).Value = 9500
$ws.Cells.Item(14,16).Value = 158

$ws.Cells.Item(15,4).Value = 44320
$ws.Cells.Item(15,11).Value = 9000
$ws.Cells.Item(15,12).Value = 10000
$ws.Cells.Item(15,13).Value = 9500
$ws.Cells.Item(15,16).Value = 158

$ws.Cells.Item(16,4).Value = 44320
$ws.Cells.Item(16,9).Value = 'Segunda'
$ws.Cells.Item(16,10).Value = 50
$ws.Cells.Item(16,11).Value = 8000
$ws.Cells.Item(16,12).Value = 8000
$ws.Cells.Item(16,13).Value = 8000
$ws.Cells.Item(16,16).Value = 133

$ws.Cells.Item(17,4).Value = 44250
$ws.Cells.Item(17,11).Value = 6000
$ws.Cells.Item(17,12).Value = 7000
$ws.Cells.Item(17,13).Value = 6500
$ws.Cells.Item(17,16).Value = 108

$ws.Cells.Item(18,4).Value = 44334
$ws.Cells.Item(18,11).Value = 11000
$ws.Cells.Item(18,12).Value = 12000
$ws.Cells.Item(18,13).Value = 11500
$ws.Cells.Item(18,16).Value = 192

$ws.Cells.Item(19,4).Value = 44299

$ws.Cells.Item(20,4).Value = 44370
$ws.Cells.Item(20,9).Value = 'Primera'
$ws.Cells.Item(20,10).Value = 100
$ws.Cells.Item(20,11).Value = 12000
$ws.Cells.Item(20,12).Value = 13000
$ws.Cells.Item(20,13).Value = 12500
$ws.Cells.Item(20,16).Value = 208

$ws.Cells.Item(21,4).Value = 44267
$ws.Cells.Item(21,11).Value = 7000
$ws.Cells.Item(21,12).Value = 8000
$ws.Cells.Item(21,13).Value = 7500
$ws.Cells.Item(21,16).Value = 125

$ws.Cells.Item(22,4).Value = 44259
$ws.Cells.Item(22,11).Value = 8000
$ws.Cells.Item(22,12).Value = 8500
$ws.Cells.Item(22,13).Value = 8250
$ws.Cells.Item(22,16).Value = 138

$ws.Cells.Item(23,4).Value = 44362
$ws.Cells.Item(23,11).Value = 11000
$ws.Cells.Item(23,12).Value = 12000
$ws.Cells.Item(23,13).Value = 11500
$ws.Cells.Item(23,16).Value = 192

$ws.Cells.Item(24,4).Value = 44433
$ws.Cells.Item(24,10).Value = 100
$ws.Cells.Item(24,11).Value = 9000
$ws.Cells.Item(24,12).Value = 10000
$ws.Cells.Item(24,13).Value = 9500
$ws.Cells.Item(24,16).Value = 158

$ws.Cells.Item(25,4).Value = 44426
$ws.Cells.Item(25,11).Value = 10000
$ws.Cells.Item(25,12).Value = 11000
$ws.Cells.Item(25,13).Value = 10500
$ws.Cells.Item(25,16).Value = 175

$ws.Cells.Item(26,4).Value = 44194
$ws.Cells.Item(26,11).Value = 7500
$ws.Cells.Item(26,12).Value = 8000
$ws.Cells.Item(26,13).Value = 7750
$ws.Cells.Item(26,16).Value = 129

$ws.Cells.Item(27,4).Value = 44398
$ws.Cells.Item(27,11).Value = 13000
$ws.Cells.Item(27,12).Value = 14000
$ws.Cells.Item(27,13).Value = 13500
$ws.Cells.Item(27,15).Value = 'Región Metropolitana'
$ws.Cells.Item(27,16).Value = 225

$ws.Cells.Item(28,4).Value = 44348
$ws.Cells.Item(28,11).Value = 12000
$ws.Cells.Item(28,12).Value = 13000
$ws.Cells.Item(28,13).Value = 12500
$ws.Cells.Item(28,16).Value = 208

$ws.Cells.Item(30,4).Value = 44341
$ws.Cells.Item(30,11).Value = 11000
$ws.Cells.Item(30,12).Value = 12000
$ws.Cells.Item(30,13).Value = 11500
$ws.Cells.Item(30,16).Value = 192

$ws.Cells.Item(31,4).Value = 44441
$ws.Cells.Item(31,11).Value = 9000
$ws.Cells.Item(31,12).Value = 10000
$ws.Cells.Item(31,13).Value = 9500
$ws.Cells.Item(31,16).Value = 158

$ws.Cells.Item(32,4).Value = 44294
$ws.Cells.Item(32,11).Value = 8000
$ws.Cells.Item(32,12).Value = 9000
$ws.Cells.Item(32,13).Value = 8500
$ws.Cells.Item(32,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(32,16).Value = 142

$ws.Cells.Item(33,4).Value = 44313
$ws.Cells.Item(33,11).Value = 8000
$ws.Cells.Item(33,12).Value = 9000
$ws.Cells.Item(33,13).Value = 8500
$ws.Cells.Item(33,16).Value = 142

$ws.Cells.Item(34,4).Value = 44435
$ws.Cells.Item(34,11).Value = 9000
$ws.Cells.Item(34,12).Value = 10000
$ws.Cells.Item(34,13).Value = 9500
$ws.Cells.Item(34,16).Value = 158

$ws.Cells.Item(35,4).Value = 44223
$ws.Cells.Item(35,11).Value = 8000
$ws.Cells.Item(35,12).Value = 8500
$ws.Cells.Item(35,13).Value = 8250
$ws.Cells.Item(35,16).Value = 138

# --- Add new row 36 (new weekly observation) ---
$ws.Cells.Item(36,1).Value = 11
$ws.Cells.Item(36,2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(36,3).Value = 'Bíobío'
$ws.Cells.Item(36,4).Value = 44238
$ws.Cells.Item(36,5).Value = 8
$ws.Cells.Item(36,6).Value = 100112001
$ws.Cells.Item(36,7).Value = 'Berenjena'
$ws.Cells.Item(36,8).Value = 'Sin especificar'
$ws.Cells.Item(36,9).Value = 'Primera'
$ws.Cells.Item(36,10).Value = 100
$ws.Cells.Item(36,11).Value = 7000
$ws.Cells.Item(36,12).Value = 8000
$ws.Cells.Item(36,13).Value = 7500
$ws.Cells.Item(36,14).Value = '$/caja 60 unidades'
$ws.Cells.Item(36,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(36,16).Value = 125
$ws.Cells.Item(36,17).Value = 60
$ws.Cells.Item(36,18).Value = 'Hortaliza'
$ws.Cells.Item(36,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

